# Apply the updated crypto price/volume snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.688.92"
$ws.Range("E2").Value = "  -0.34%  "

# Row 3
$ws.Range("D3").Value = "1.603.16"
$ws.Range("E3").Value = "  -0.40%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("E5").Value = "  -0.54%  "

# Row 6
$ws.Range("E6").Value = "  +0.43%  "

# Row 7
$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = $ws.Range("D6").Style
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("D8").Value = "'29.42"
$ws.Range("D8").Style = $ws.Range("D9").Style
$ws.Range("E8").Value = "  +9.24%  "

# Row 9
$ws.Range("E9").Value = "  +2.34%  "

# Row 10
$ws.Range("E10").Value = "  +1.18%  "

# Row 11
$ws.Range("D11").Value = "'0.0904"
$ws.Range("D11").Style = $ws.Range("D10").Style
$ws.Range("E11").Value = "  -0.87%  "

# Row 12
$ws.Range("D12").Value = "1.832.92"
$ws.Range("E12").Value = "  -0.40%  "

# Row 13
$ws.Range("D13").Value = "1.603.25"
$ws.Range("E13").Value = "  -0.33%  "

# Row 14
$ws.Range("D14").Value = "'0.555"
$ws.Range("D14").Style = $ws.Range("D10").Style
$ws.Range("E14").Value = "  +3.00%  "

# Row 15
$ws.Range("D15").Value = "29.708.64"
$ws.Range("E15").Value = "  -0.37%  "

# Row 16
$ws.Range("D16").Value = "'3.80"
$ws.Range("D16").Style = $ws.Range("D21").Style
$ws.Range("E16").Value = "  +1.06%  "

# Row 17
$ws.Range("D17").Value = "'64.15"
$ws.Range("D17").Style = $ws.Range("D21").Style
$ws.Range("E17").Value = "  +1.27%  "

# Row 18
$ws.Range("D18").Value = "'242.71"
$ws.Range("D18").Style = $ws.Range("D21").Style
$ws.Range("E18").Value = "  -0.86%  "

# Row 19
$ws.Range("D19").Value = "'8.01"
$ws.Range("D19").Style = $ws.Range("D21").Style
$ws.Range("E19").Value = "  +5.38%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0700"
$ws.Range("E20").Value = "  +0.91%  "

# Row 21
$ws.Range("E21").Value = "  +0.14%  "

# Row 22
$ws.Range("E22").Value = "  -0.05%  "

# Row 23
$ws.Range("D23").Value = "'9.53"
$ws.Range("D23").Style = $ws.Range("D22").Style
$ws.Range("E23").Value = "  +3.05%  "

# Row 24
$ws.Range("D24").Value = "'2.11"
$ws.Range("D24").Style = $ws.Range("D22").Style
$ws.Range("E24").Value = "  +0.77%  "

# Row 25
$ws.Range("D25").Value = "'155.50"
$ws.Range("D25").Style = $ws.Range("D27").Style
$ws.Range("E25").Value = "  -0.30%  "

# Row 26
$ws.Range("D26").Value = "'15.54"
$ws.Range("D26").Style = $ws.Range("D27").Style
$ws.Range("E26").Value = "  +1.35%  "

# Row 27
$ws.Range("E27").Value = "  +0.88%  "

# Row 28
$ws.Range("D28").Value = "'6.47"
$ws.Range("D28").Style = $ws.Range("D27").Style
$ws.Range("E28").Value = "  +1.33%  "

# Row 29
$ws.Range("D29").Value = "'0.997"
$ws.Range("D29").Style = $ws.Range("D31").Style
$ws.Range("E29").Value = "  +0.09%  "

# Row 30
$ws.Range("D30").Value = "'0.0480"
$ws.Range("D30").Style = $ws.Range("D31").Style
$ws.Range("E30").Value = "  +1.50%  "

# Row 31
$ws.Range("E31").Value = "  -0.13%  "

# Row 32
$ws.Range("E32").Value = "  -0.15%  "

# Row 33
$ws.Range("E33").Value = "  +2.39%  "

# Row 34
$ws.Range("D34").Value = "1.426.38"
$ws.Range("E34").Value = "  -1.42%  "

# Row 35
$ws.Range("D35").Value = "'1.58"
$ws.Range("D35").Style = $ws.Range("D36").Style
$ws.Range("E35").Value = "  +4.40%  "

# Row 36
$ws.Range("E36").Value = "  -0.23%  "

# Row 37
$ws.Range("E37").Value = "  +1.04%  "

# Row 38
$ws.Range("E38").Value = "  +0.48%  "

# Row 39
$ws.Range("E39").Value = "  +1.76%  "

# Row 40
$ws.Range("D40").Value = "'0.547"
$ws.Range("D40").Style = $ws.Range("D41").Style
$ws.Range("E40").Value = "  +2.40%  "

# Row 41
$ws.Range("E41").Value = "  +0.29%  "

# Row 42
$ws.Range("E42").Value = "  +3.32%  "

# Row 43
$ws.Range("D43").Value = "'0.0494"
$ws.Range("D43").Style = $ws.Range("D44").Style
$ws.Range("E43").Value = "  +5.57%  "

# Row 44
$ws.Range("E44").Value = "  +0.77%  "

# Row 45
$ws.Range("E45").Value = "  +0.11%  "

# Row 47
$ws.Range("E47").Value = "  +1.37%  "

# Row 48
$ws.Range("D48").Value = "'5.38"
$ws.Range("D48").Style = $ws.Range("D47").Style
$ws.Range("E48").Value = "  +1.38%  "

# Row 49
$ws.Range("D49").Value = "1.742.39"
$ws.Range("E49").Value = "  -0.59%  "

# Row 50
$ws.Range("D50").Value = "'86.84"
$ws.Range("D50").Style = $ws.Range("D47").Style
$ws.Range("E50").Value = "  -0.09%  "

# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0104"
$ws.Range("E51").Value = "  +2.40%  "
